# Apply the commit's changes:
#  - Sheet "Y": add two new columns (F: "qS/n" / "inp3", G: "qS/qS" / 166)
#    and make it the active/selected sheet with G5 selected.
#  - Sheet "M": no longer the active sheet; selection moves to F11.

$wb = $excel.ActiveWorkbook

$wsY = $wb.Worksheets.Item("Y")
$wsM = $wb.Worksheets.Item("M")

# New header + data cells in columns F and G of sheet "Y"
$wsY.Cells.Item(1, 6).Value = "qS/n"
$wsY.Cells.Item(1, 7).Value = "qS/qS"
$wsY.Cells.Item(2, 6).Value = "inp3"
$wsY.Cells.Item(2, 7).Value = 166

# Sheet "M" keeps a (new) selection even though it's no longer the active tab.
$wsM.Range("F11").Select()

# Sheet "Y" becomes the active sheet, with G5 selected.
$wsY.Activate()
$wsY.Range("G5").Select()
